# Update the "RecipientTo" (To List) value in the historical/archived
# ConfigOptions blocks so they match the current Config sheet's To list
# (adds david.villasoto@lexisnexisrisk.com to the recipients).

$wb = $excel.ActiveWorkbook

$newTo = "sam.tecson@lexisnexisrisk.com; joavic.quisano@lexisnexisrisk.com; david.villasoto@lexisnexisrisk.com"

$configOptions = $wb.Worksheets.Item("ConfigOptions")
$configOptions.Range("B11").Value = $newTo
$configOptions.Range("B25").Value = $newTo
$configOptions.Range("B39").Value = $newTo

# Restore selections/active cells to match what was left selected when
# the workbook was saved.
$config = $wb.Worksheets.Item("Config")

$configOptions.Activate()
$configOptions.Range("B39:B40").Select()

$config.Activate()
$config.Range("B19").Select()
